# feat: add 2022-Q1 data
#
# The existing "总计" (Total) sheet's position becomes the new "2022-Q1" sheet
# (fund holdings for the new quarter), and a brand-new "总计" sheet is appended
# after it, containing the updated totals table (with a new 2022-Q1 row on top,
# old rows shifted down by one).
#
# NOTE: worksheet object variables can become stale after a structural change
# (Copy/Insert shifts what a previously-captured reference resolves to), so we
# re-fetch sheets by name right before using them whenever the sheet list has
# just changed.

$wb = $excel.ActiveWorkbook

# --- Step 1: copy the "2021-Q4" sheet (same fund-holdings layout, 3 data rows)
#     to build the new "2022-Q1" sheet, inserting it right before "总计". ---
$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$prevQuarter.Copy($totalSheet, $null)

$newQ1 = $wb.Worksheets.Item("2021-Q4 (2)")
$newQ1.Name = "2022-Q1"

# Overwrite the fund-holding rows with the 2022-Q1 data (header row stays the same).
#
# Columns B/D/E/F/G hold digit-looking text ("000880", "3.98", ...) that must
# stay plain text (no leading-zero/precision loss). Forcing NumberFormat="@"
# before assigning the value keeps it text, but leaves a permanent "Text"
# style stamped on the cell; the source sheets don't have any such style on
# these cells at all. So after writing the value, reset the cell's format by
# pasting-special (formats only) from a guaranteed-blank cell -- that clears
# the stamped style index back to "no style", matching the source layout,
# while leaving the text value itself untouched.
$blank = $newQ1.Cells.Item(50, 50)

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $blank.Copy()
    $cell.PasteSpecial(-4122)
}

Set-TextValue $newQ1.Cells.Item(2, 2) "000880"
$newQ1.Cells.Item(2, 3).Value = "富国研究精选灵活配置混合"
Set-TextValue $newQ1.Cells.Item(2, 4) "3.98"
Set-TextValue $newQ1.Cells.Item(2, 5) "90.99"
Set-TextValue $newQ1.Cells.Item(2, 6) "7.90"
Set-TextValue $newQ1.Cells.Item(2, 7) "0.3144"
$newQ1.Cells.Item(2, 8).Value = 3

Set-TextValue $newQ1.Cells.Item(3, 2) "008736"
$newQ1.Cells.Item(3, 3).Value = "南方高股息主题股票A"
Set-TextValue $newQ1.Cells.Item(3, 4) "1.02"
Set-TextValue $newQ1.Cells.Item(3, 5) "89.27"
Set-TextValue $newQ1.Cells.Item(3, 6) "4.23"
Set-TextValue $newQ1.Cells.Item(3, 7) "0.0431"
$newQ1.Cells.Item(3, 8).Value = 5

Set-TextValue $newQ1.Cells.Item(4, 2) "008737"
$newQ1.Cells.Item(4, 3).Value = "南方高股息主题股票C"
Set-TextValue $newQ1.Cells.Item(4, 4) "0.08"
Set-TextValue $newQ1.Cells.Item(4, 5) "89.27"
Set-TextValue $newQ1.Cells.Item(4, 6) "4.23"
Set-TextValue $newQ1.Cells.Item(4, 7) "0.0034"
$newQ1.Cells.Item(4, 8).Value = 5

$blank.Clear() | Out-Null

# --- Step 2: add a new 2022-Q1 row at the top of the real "总计" sheet's data,
#     shifting the existing rows down by one. Re-fetch "总计" fresh since the
#     sheet list just changed above. Written as a direct rewrite of the final
#     table (instead of Rows.Insert, which drags unwanted formatting onto the
#     shifted cells) so every cell ends up with the same styling as before. ---
$realTotal = $wb.Worksheets.Item("总计")

# Stamp row 7's index cell with the same style as the other index cells (A2:A6)
# before filling it in, since it's a brand-new row.
$realTotal.Range("A2").Copy()
$realTotal.Range("A7").PasteSpecial(-4122)

$realTotal.Cells.Item(7, 1).Value = 5
$realTotal.Cells.Item(7, 2).Value = "2020-Q4"
$realTotal.Cells.Item(7, 3).Value = 6
$realTotal.Cells.Item(7, 4).Value = 0.76

$realTotal.Cells.Item(6, 1).Value = 4
$realTotal.Cells.Item(6, 2).Value = "2021-Q1"
$realTotal.Cells.Item(6, 3).Value = 12
$realTotal.Cells.Item(6, 4).Value = 0.97

$realTotal.Cells.Item(5, 1).Value = 3
$realTotal.Cells.Item(5, 2).Value = "2021-Q2"
$realTotal.Cells.Item(5, 3).Value = 9
$realTotal.Cells.Item(5, 4).Value = 0.7

$realTotal.Cells.Item(4, 1).Value = 2
$realTotal.Cells.Item(4, 2).Value = "2021-Q3"
$realTotal.Cells.Item(4, 3).Value = 5
$realTotal.Cells.Item(4, 4).Value = 0.38

$realTotal.Cells.Item(3, 1).Value = 1
$realTotal.Cells.Item(3, 2).Value = "2021-Q4"
$realTotal.Cells.Item(3, 3).Value = 3
$realTotal.Cells.Item(3, 4).Value = 4.34

$realTotal.Cells.Item(2, 1).Value = 0
$realTotal.Cells.Item(2, 2).Value = "2022-Q1"
$realTotal.Cells.Item(2, 3).Value = 3
$realTotal.Cells.Item(2, 4).Value = 0.36
